$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# This edit reorders several player rows (A:C) within the sheet.
# The underlying permutation decomposes into:
#   - swap(3, 14)
#   - rotate(4 -> gets old 5, 5 -> gets old 15, 15 -> gets old 4)
#   - swap(11, 16)

function Get-RowValues($row) {
    return @(
        $ws.Range("A$row").Value2,
        $ws.Range("B$row").Value2,
        $ws.Range("C$row").Value2
    )
}

function Set-RowValues($row, $vals) {
    $ws.Range("A$row").Value = $vals[0]
    $ws.Range("B$row").Value = $vals[1]
    $ws.Range("C$row").Value = $vals[2]
}

# Capture original values for all affected rows before any writes.
$row3  = Get-RowValues 3
$row4  = Get-RowValues 4
$row5  = Get-RowValues 5
$row11 = Get-RowValues 11
$row14 = Get-RowValues 14
$row15 = Get-RowValues 15
$row16 = Get-RowValues 16

# Swap rows 3 and 14
Set-RowValues 3  $row14
Set-RowValues 14 $row3

# Rotate rows 4, 5, 15: 4<-5, 5<-15, 15<-4
Set-RowValues 4  $row5
Set-RowValues 5  $row15
Set-RowValues 15 $row4

# Swap rows 11 and 16
Set-RowValues 11 $row16
Set-RowValues 16 $row11
